{"js": "// Mobile acceptance test read-me: proofreading pass.\n//  - \"Umbuy\" (first mention, bold heading)      -> \"UMBUY\"\n//  - \"allow\" (bold heading, verb)               -> \"allowed\"\n//  - \"UmBuy\" (Settings navigation instruction)   -> \"UMBUY\"\n//  - move the \"_GoBack\" bookmark from the trailing empty paragraph to\n//    right after the \"...UMBUY\" run in the Settings navigation line\n//    (this is where Word leaves it after the last edit before save).\n\nconst body = context.document.body;\n\n// 1. \"Tester should be logged out ... Acceptance Test for Umbuy mobile\n//    application.\" -> capitalize the product name.\nconst umbuy1 = body.search(\"Umbuy\", { matchCase: true, matchWholeWord: true });\numbuy1.load(\"items\");\nawait context.sync();\numbuy1.items[0].insertText(\"UMBUY\", \"Replace\");\nawait context.sync();\n\n// 2. \"Acceptance Test doesn't cover if the user hasn't allow the mobile\n//    application ...\" -> \"... hasn't allowed the mobile application ...\"\nconst allow = body.search(\"allow\", { matchCase: true, matchWholeWord: true });\nallow.load(\"items\");\nawait context.sync();\nallow.items[0].insertText(\"allowed\", \"Replace\");\nawait context.sync();\n\n// 3. \"Tester should go to the Settings -> Apps & Notifications -> UmBuy ->\n//    Permissions -> Give permission to Storage.\" -> capitalize the app name.\nconst umbuy2 = body.search(\"UmBuy\", { matchCase: true, matchWholeWord: true });\numbuy2.load(\"items\");\nawait context.sync();\numbuy2.items[0].insertText(\"UMBUY\", \"Replace\");\nawait context.sync();\n\n// 4. Relocate the \"_GoBack\" bookmark (Word keeps only one instance of it,\n//    marking the last edit position) from the final empty paragraph to\n//    just after the \"...UMBUY\" run we edited above.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst umbuyRuns = body.search(\"UMBUY\", { matchCase: true, matchWholeWord: true });\numbuyRuns.load(\"items\");\nawait context.sync();\nconst settingsUmbuy = umbuyRuns.items[umbuyRuns.items.length - 1];\nconst collapsedEnd = settingsUmbuy.getRange(\"End\");\ncollapsedEnd.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Mobile acceptance test read-me: proofreading pass.\n#  - \"Umbuy\" (first mention, bold heading)      -> \"UMBUY\"\n#  - \"allow\" (bold heading, verb)               -> \"allowed\"\n#  - \"UmBuy\" (Settings navigation instruction)   -> \"UMBUY\"\n#  - move the \"_GoBack\" bookmark from the trailing empty paragraph to\n#    right after the \"...UMBUY\" run in the Settings navigation line\n#    (this is where Word leaves it after the last edit before save).\n\n$d = $word.ActiveDocument\n\n# 1. \"Tester should be logged out ... Acceptance Test for Umbuy mobile\n#    application.\" -> capitalize the product name.\n$d.Content.Find.Execute(\"Umbuy\", $true, $true, $false, $false, $false, $true, 1, $false, \"UMBUY\", 2)\n\n# 2. \"Acceptance Test doesn't cover if the user hasn't allow the mobile\n#    application ...\" -> \"... hasn't allowed the mobile application ...\"\n$d.Content.Find.Execute(\"allow\", $true, $true, $false, $false, $false, $true, 1, $false, \"allowed\", 2)\n\n# 3. \"Tester should go to the Settings -> Apps & Notifications -> UmBuy ->\n#    Permissions -> Give permission to Storage.\" -> capitalize the app name.\n$d.Content.Find.Execute(\"UmBuy\", $true, $true, $false, $false, $false, $true, 1, $false, \"UMBUY\", 2)\n\n# 4. Relocate the \"_GoBack\" bookmark (Word keeps only one instance of it,\n#    marking the last edit position) from the final empty paragraph to\n#    just after the \"...UMBUY\" run we edited in step 3 above.\n#    Find the last \"UMBUY\" occurrence in the document (the one from step 3).\n$docEnd = $d.Content.End\n$searchStart = 0\n$lastMatchEnd = -1\nwhile ($true) {\n    $rng = $d.Range($searchStart, $docEnd)\n    $found = $rng.Find.Execute(\"UMBUY\", $true, $true, $false, $false, $false, $true, 1, $false)\n    if (-not $found) { break }\n    $lastMatchEnd = $rng.End\n    $searchStart = $rng.End\n}\n\n$d.Bookmarks(\"_GoBack\").Delete()\n\n$bmRange = $d.Range($lastMatchEnd, $lastMatchEnd)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange)\n"}
